$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-13 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-14 Thursday", 2) | Out-Null
$d.Content.Find.Execute("23×51=", $true, $false, $false, $false, $false, $true, 1, $false, "12×29=", 2) | Out-Null
$d.Content.Find.Execute("11×93=", $true, $false, $false, $false, $false, $true, 1, $false, "55×64=", 2) | Out-Null
$d.Content.Find.Execute("25×67=", $true, $false, $false, $false, $false, $true, 1, $false, "52×69=", 2) | Out-Null
$d.Content.Find.Execute("53×86=", $true, $false, $false, $false, $false, $true, 1, $false, "27×29=", 2) | Out-Null
$d.Content.Find.Execute("15×32=", $true, $false, $false, $false, $false, $true, 1, $false, "84×79=", 2) | Out-Null
$d.Content.Find.Execute("13×26=", $true, $false, $false, $false, $false, $true, 1, $false, "23×93=", 2) | Out-Null
$d.Content.Find.Execute("72×46=", $true, $false, $false, $false, $false, $true, 1, $false, "71×42=", 2) | Out-Null
$d.Content.Find.Execute("71×73=", $true, $false, $false, $false, $false, $true, 1, $false, "47×33=", 2) | Out-Null
$d.Content.Find.Execute("82×91=", $true, $false, $false, $false, $false, $true, 1, $false, "84×87=", 2) | Out-Null
$d.Content.Find.Execute("76×50=", $true, $false, $false, $false, $false, $true, 1, $false, "26×89=", 2) | Out-Null
$d.Content.Find.Execute("46×61=", $true, $false, $false, $false, $false, $true, 1, $false, "13×60=", 2) | Out-Null
$d.Content.Find.Execute("19×80=", $true, $false, $false, $false, $false, $true, 1, $false, "28×44=", 2) | Out-Null
$d.Content.Find.Execute("34×38=", $true, $false, $false, $false, $false, $true, 1, $false, "47×64=", 2) | Out-Null
$d.Content.Find.Execute("34×32=", $true, $false, $false, $false, $false, $true, 1, $false, "42×64=", 2) | Out-Null
$d.Content.Find.Execute("45×21=", $true, $false, $false, $false, $false, $true, 1, $false, "97×99=", 2) | Out-Null
$d.Content.Find.Execute("12×97=", $true, $false, $false, $false, $false, $true, 1, $false, "47×87=", 2) | Out-Null
$d.Content.Find.Execute("83×67=", $true, $false, $false, $false, $false, $true, 1, $false, "85×65=", 2) | Out-Null
$d.Content.Find.Execute("97×50=", $true, $false, $false, $false, $false, $true, 1, $false, "39×17=", 2) | Out-Null
$d.Content.Find.Execute("95×76=", $true, $false, $false, $false, $false, $true, 1, $false, "99×24=", 2) | Out-Null
$d.Content.Find.Execute("42×55=", $true, $false, $false, $false, $false, $true, 1, $false, "50×53=", 2) | Out-Null
$d.Content.Find.Execute("15×21=", $true, $false, $false, $false, $false, $true, 1, $false, "39×42=", 2) | Out-Null
$d.Content.Find.Execute("54×58=", $true, $false, $false, $false, $false, $true, 1, $false, "45×13=", 2) | Out-Null
$d.Content.Find.Execute("94×80=", $true, $false, $false, $false, $false, $true, 1, $false, "41×85=", 2) | Out-Null
$d.Content.Find.Execute("81×13=", $true, $false, $false, $false, $false, $true, 1, $false, "71×57=", 2) | Out-Null
$d.Content.Find.Execute("55×98=", $true, $false, $false, $false, $false, $true, 1, $false, "69×80=", 2) | Out-Null
